$wb = $excel.ActiveWorkbook
$wsGames = $wb.Worksheets.Item("Games")
$wsNext = $wb.Worksheets.Item("Next")

# --- "Next" sheet: the scheduled game (row 2: 2024-01-15 @ CLE) has now
# been played, so it drops off the top of the schedule and everything
# else shifts up one row.
$wsNext.Rows.Item(2).Delete()

# --- "Games" sheet: append the completed game (the one that dropped off
# the "Next" schedule) with its final boxscore stats as the new last row.
$newRow = $wsGames.UsedRange.Rows.Count + 1

$wsGames.Range("A$newRow").Value = 42
$wsGames.Range("B$newRow").Value = 45306
$wsGames.Range("B$newRow").NumberFormat = "YYYY-MM-DD"
$wsGames.Range("C$newRow").Value = -1
$wsGames.Range("D$newRow").Value = 91
$wsGames.Range("E$newRow").Value = 97.40000000000001
$wsGames.Range("F$newRow").Value = 0.476
$wsGames.Range("G$newRow").Value = 16.7
$wsGames.Range("H$newRow").Value = 14.6
$wsGames.Range("I$newRow").Value = 0.145
$wsGames.Range("J$newRow").Value = 93.5
$wsGames.Range("K$newRow").Value = "CLE"
$wsGames.Range("L$newRow").Value = 109
$wsGames.Range("M$newRow").Value = 0.516
$wsGames.Range("N$newRow").Value = 9.9
$wsGames.Range("O$newRow").Value = 26
$wsGames.Range("P$newRow").Value = 0.116
$wsGames.Range("Q$newRow").Value = 112
$wsGames.Range("R$newRow").Value = 0
$wsGames.Range("S$newRow").Value = 0
